$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.087.77"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "2.369.10"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").Value = "2.735.71"
$ws.Range("E15").Value = "  +1.14%  "

$ws.Range("D16").Value = "2.372.06"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.803"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "43.104.91"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("E29").Value = "  +1.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0740"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "128.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.17%  "

$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.70%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.75%  "

$ws.Range("D43").Value = "1.929.68"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.66%  "

$ws.Range("D48").Value = "2.596.20"
$ws.Range("E48").Value = "  +1.00%  "

$ws.Range("E49").Value = "  +3.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "
